$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.074334117552723
$ws.Range("C2").Value = 0.2111556234066256
$ws.Range("D2").Value = 0.07589854406890595
$ws.Range("E2").Value = 0.02972752957218194
$ws.Range("G2").Value = 0.002646906146153501
$ws.Range("I2").Value = 5.689741979547478
$ws.Range("K2").Value = 1.168990046916491
$ws.Range("L2").Value = 0.2975616675253576

$ws.Range("B3").Value = 1.073948960523239
$ws.Range("C3").Value = 0.1924337519348285
$ws.Range("D3").Value = 0.06924553444146397
$ws.Range("E3").Value = 0.0299793619647204
$ws.Range("G3").Value = 0.002652928135165718
$ws.Range("I3").Value = 5.359892256295893
$ws.Range("K3").Value = 1.145283538329238
$ws.Range("L3").Value = 0.2890363637149989

$ws.Range("B4").Value = 1.075003360951683
$ws.Range("C4").Value = 0.1809936237831096
$ws.Range("D4").Value = 0.06520433084868671
$ws.Range("E4").Value = 0.0301527770937895
$ws.Range("G4").Value = 0.002656811085476234
$ws.Range("I4").Value = 5.156806094405141
$ws.Range("K4").Value = 1.131886134249783
$ws.Range("L4").Value = 0.2839698418084851

$ws.Range("B5").Value = 1.075756309818189
$ws.Range("C5").Value = 0.1763451283723043
$ws.Range("D5").Value = 0.06356832596226525
$ws.Range("E5").Value = 0.03022817265433453
$ws.Range("G5").Value = 0.002658440228613066
$ws.Range("I5").Value = 5.073897132045857
$ws.Range("K5").Value = 1.126716361719019
$ws.Range("L5").Value = 0.2819472387981676

$ws.Range("B6").Value = 1.075900822848837
$ws.Range("C6").Value = 0.1755740493134965
$ws.Range("D6").Value = 0.06329731656492754
$ws.Range("E6").Value = 0.03024097764560629
$ws.Range("G6").Value = 0.002658713579109074
$ws.Range("I6").Value = 5.060120841756742
$ws.Range("K6").Value = 1.125875382678601
$ws.Range("L6").Value = 0.2816139205104093

$ws.Range("B7").Value = 1.07501220829559
$ws.Range("C7").Value = 0.1809308786476436
$ws.Range("D7").Value = 0.06518222355741443
$ws.Range("E7").Value = 0.03015377475919045
$ws.Range("G7").Value = 0.002656832866867009
$ws.Range("I7").Value = 5.155688573922873
$ws.Range("K7").Value = 1.131815241580512
$ws.Range("L7").Value = 0.2839423943088946

$ws.Range("B8").Value = 1.073932575024543
$ws.Range("C8").Value = 0.2046886068584968
$ws.Range("D8").Value = 0.0735953853725988
$ws.Range("E8").Value = 0.02981046333363402
$ws.Range("G8").Value = 0.002648944153322797
$ws.Range("I8").Value = 5.576119509886297
$ws.Range("K8").Value = 1.160574654073514
$ws.Range("L8").Value = 0.2945871319302995

$ws.Range("B9").Value = 1.082122153093621
$ws.Range("C9").Value = 0.2517369893205057
$ws.Range("D9").Value = 0.09045034988771761
$ws.Range("E9").Value = 0.02928621406136234
$ws.Range("G9").Value = 0.002634937299056481
$ws.Range("I9").Value = 6.396636944722815
$ws.Range("K9").Value = 1.226239681180289
$ws.Range("L9").Value = 0.3168060165587718

$ws.Range("B10").Value = 1.094516979205508
$ws.Range("C10").Value = 0.2866174384332396
$ws.Range("D10").Value = 0.1030657395268406
$ws.Range("E10").Value = 0.02899177779832307
$ws.Range("G10").Value = 0.002625526496638023
$ws.Range("I10").Value = 6.997779804645063
$ws.Range("K10").Value = 1.280248470232834
$ws.Range("L10").Value = 0.3339678213670823

$ws.Range("B11").Value = 1.101562023745402
$ws.Range("C11").Value = 0.3025612921090044
$ws.Range("D11").Value = 0.1088584734513063
$ws.Range("E11").Value = 0.02887751924505189
$ws.Range("G11").Value = 0.00262143383854454
$ws.Range("I11").Value = 7.271047796928883
$ws.Range("K11").Value = 1.306096229573598
$ws.Range("L11").Value = 0.3419611579002719

$ws.Range("B12").Value = 1.104433798696107
$ws.Range("C12").Value = 0.3086103850343136
$ws.Range("D12").Value = 0.1110600311998127
$ws.Range("E12").Value = 0.02883708185312983
$ws.Range("G12").Value = 0.002619910946453974
$ws.Range("I12").Value = 7.374511301917721
$ws.Range("K12").Value = 1.3160700378819
$ws.Range("L12").Value = 0.3450151244876167

$ws.Range("B13").Value = 1.103806211331005
$ws.Range("C13").Value = 0.3073070851594082
$ws.Range("D13").Value = 0.1105855280981132
$ws.Range("E13").Value = 0.02884566491450258
$ws.Range("G13").Value = 0.002620237734724152
$ws.Range("I13").Value = 7.352229172527927
$ws.Range("K13").Value = 1.313913704174524
$ws.Range("L13").Value = 0.3443561901428751

$ws.Range("B14").Value = 1.101794188617873
$ws.Range("C14").Value = 0.3030587212584237
$ws.Range("D14").Value = 0.1090394358211881
$ws.Range("E14").Value = 0.02887413572358533
$ws.Range("G14").Value = 0.002621308011241811
$ws.Range("I14").Value = 7.279560082526416
$ws.Range("K14").Value = 1.306913045442997
$ws.Range("L14").Value = 0.3422118656161075

$ws.Range("B15").Value = 1.100588380845494
$ws.Range("C15").Value = 0.3004579883287022
$ws.Range("D15").Value = 0.1080934548689783
$ws.Range("E15").Value = 0.02889194343793555
$ws.Range("G15").Value = 0.002621967085311222
$ws.Range("I15").Value = 7.235046258932698
$ws.Range("K15").Value = 1.302649197289185
$ws.Range("L15").Value = 0.3409019374670947

$ws.Range("B16").Value = 1.094085017649604
$ws.Range("C16").Value = 0.2855770525452499
$ws.Range("D16").Value = 0.1026882746552076
$ws.Range("E16").Value = 0.02899964087472817
$ws.Range("G16").Value = 0.00262579773666472
$ws.Range("I16").Value = 6.979918101727009
$ws.Range("K16").Value = 1.278585155890795
$ws.Range("L16").Value = 0.3334492104384026

$ws.Range("B17").Value = 1.090456872133984
$ws.Range("C17").Value = 0.2764680366229015
$ws.Range("D17").Value = 0.09938633952573639
$ws.Range("E17").Value = 0.02907075042087826
$ws.Range("G17").Value = 0.002628195835397226
$ws.Range("I17").Value = 6.823362995561126
$ws.Range("K17").Value = 1.264151550639042
$ws.Range("L17").Value = 0.3289251189863052

$ws.Range("B18").Value = 1.088502347690593
$ws.Range("C18").Value = 0.2712359479388624
$ws.Range("D18").Value = 0.09749221646387696
$ws.Range("E18").Value = 0.02911350332034601
$ws.Range("G18").Value = 0.002629592898556632
$ws.Range("I18").Value = 6.73329763114603
$ws.Range("K18").Value = 1.255969980024844
$ws.Range("L18").Value = 0.3263405044877885

$ws.Range("B19").Value = 1.08786324878767
$ws.Range("C19").Value = 0.2694656703419867
$ws.Range("D19").Value = 0.09685176183410249
$ws.Range("E19").Value = 0.02912829691440422
$ws.Range("G19").Value = 0.002630068972152297
$ws.Range("I19").Value = 6.702799453513876
$ws.Range("K19").Value = 1.253220437417241
$ws.Range("L19").Value = 0.3254684000504824

$ws.Range("B20").Value = 1.090829390812218
$ws.Range("C20").Value = 0.277436960662186
$ws.Range("D20").Value = 0.0997373106089583
$ws.Range("E20").Value = 0.0290629889581755
$ws.Range("G20").Value = 0.002627938718870209
$ws.Range("I20").Value = 6.840030442981686
$ws.Range("K20").Value = 1.26567557263246
$ws.Range("L20").Value = 0.3294049007448194

$ws.Range("B21").Value = 1.102379618968058
$ws.Range("C21").Value = 0.304306252595353
$ws.Range("D21").Value = 0.109493342381441
$ws.Range("E21").Value = 0.0288656963576468
$ws.Range("G21").Value = 0.002620992915814746
$ws.Range("I21").Value = 7.300905133662809
$ws.Range("K21").Value = 1.308964249880262
$ws.Range("L21").Value = 0.3428409690430243

$ws.Range("B22").Value = 1.111118075006942
$ws.Range("C22").Value = 0.3219341640720472
$ws.Range("D22").Value = 0.1159160638129606
$ws.Range("E22").Value = 0.02875324872859508
$ws.Range("G22").Value = 0.002616610195524689
$ws.Range("I22").Value = 7.602018594445383
$ws.Range("K22").Value = 1.338339947734312
$ws.Range("L22").Value = 0.3517800931138737

$ws.Range("B23").Value = 1.106344762222506
$ws.Range("C23").Value = 0.3125194948581793
$ws.Range("D23").Value = 0.1124837995092349
$ws.Range("E23").Value = 0.0288117549607172
$ws.Range("G23").Value = 0.00261893505090275
$ws.Range("I23").Value = 7.441313702084358
$ws.Range("K23").Value = 1.32256174717034
$ws.Range("L23").Value = 0.3469945782198494

$ws.Range("B24").Value = 1.090660566136677
$ws.Range("C24").Value = 0.2769988951714026
$ws.Range("D24").Value = 0.09957862350509572
$ws.Range("E24").Value = 0.02906649208748213
$ws.Range("G24").Value = 0.002628054903837946
$ws.Range("I24").Value = 6.832495276925272
$ws.Range("K24").Value = 1.264986199558393
$ws.Range("L24").Value = 0.3291879405180111

$ws.Range("B25").Value = 1.078793611657915
$ws.Range("C25").Value = 0.2389564940488356
$ws.Range("D25").Value = 0.08585095188412595
$ws.Range("E25").Value = 0.02941209837292114
$ws.Range("G25").Value = 0.00263857113052295
$ws.Range("I25").Value = 6.175014519750476
$ws.Range("K25").Value = 1.207471667211138
$ws.Range("L25").Value = 0.3106494496813212
